# "completed desktop recommendation part"
#
# The Desktop recommendation sheet gets a blank leading row (everything
# shifts down by one), plus a separator row containing a single "!" cell
# is inserted before each of the three lower product blocks. The third
# block ("AMD Ryzen 7 5800X" build) also gains a graphics-card line
# (MSI GeForce RTX 3060, $650) that folds into its Total formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole sheet down by one row: old row 1 -> row 2, row 2 -> row 3,
# row 4 -> row 5, row 5 -> row 6, row 7 -> row 8, row 8 -> row 9,
# row 10 -> row 11, row 11 -> row 12, row 12 -> row 13, row 13 -> row 14.
# This also leaves rows 4, 7 and 10 empty - exactly where the new "!"
# separator rows belong - and Excel auto-adjusts the SUM formulas in
# column H for us.
$ws.Rows("1:1").Insert()

# New "!" separator rows between the product blocks.
$ws.Cells.Item(4, 1).Value2 = "!"
$ws.Cells.Item(7, 1).Value2 = "!"
$ws.Cells.Item(10, 1).Value2 = "!"

# New graphics-card row for the third build (now at sheet rows 8/9).
$ws.Cells.Item(8, 7).Value2 = "MSI GeForce RTX 3060"
$ws.Cells.Item(9, 7).Value2 = 650
$ws.Cells.Item(9, 8).Formula = "=SUM(A9:G9)"

# Restore the active-cell selection recorded in the saved workbook.
$ws.Range("E8").Select() | Out-Null
